$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 221.9
$ws.Range("I2").Value = 193.22223
$ws.Range("K2").Value = 193.22223
$ws.Range("M2").Value = -80.22223
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
$ws.Range("H6").Value = 1322.8518
$ws.Range("I6").Value = 1256.6842
$ws.Range("K6").Value = 3770.0526
$ws.Range("M6").Value = -3658.0526
$ws.Range("H33").Value = 415.30768
$ws.Range("I33").Value = 208.81818
$ws.Range("J33").Value = 1551
$ws.Range("K33").Value = 208.81818
$ws.Range("L33").Value = 1551
$ws.Range("M33").Value = 20.18181999999999
$ws.Range("N33").Value = -2009
$ws.Range("H41").Value = 149.26666
$ws.Range("J41").Value = 134
$ws.Range("L41").Value = 134
$ws.Range("N41").Value = -1014
$ws.Range("H98").Value = 5027.5557
$ws.Range("I98").Value = 5027.5557
$ws.Range("K98").Value = 5027.5557
$ws.Range("M98").Value = -3529.5557
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H122").Value = 5027.5557
$ws.Range("I122").Value = 5027.5557
$ws.Range("K122").Value = 15082.6671
$ws.Range("M122").Value = -12632.6671
$ws.Range("H125").Value = 7964.6206
$ws.Range("J125").Value = 7198.9375
$ws.Range("L125").Value = 64790.4375
$ws.Range("N125").Value = -69710.4375
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
$ws.Range("H135").Value = 1812
$ws.Range("I135").Value = 1416
$ws.Range("K135").Value = 12744
$ws.Range("M135").Value = -10209
$ws.Range("H137").Value = 2007.0588
$ws.Range("J137").Value = 2679.7144
$ws.Range("L137").Value = 8039.1432
$ws.Range("N137").Value = -13139.1432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2850.077
$ws.Range("I2").Value = 2196.5789
$ws.Range("J2").Value = 4623.857
$ws.Range("K2").Value = 2196.5789
$ws.Range("L2").Value = 4623.857
$ws.Range("M2").Value = -2083.5789
$ws.Range("N2").Value = -4849.857
$ws.Range("H32").Value = 4753.4443
$ws.Range("I32").Value = 4731.017
$ws.Range("K32").Value = 4731.017
$ws.Range("M32").Value = -4444.017
$ws.Range("H61").Value = 19236778
$ws.Range("I61").Value = 22732012
$ws.Range("K61").Value = 22732012
$ws.Range("M61").Value = -22731800
$ws.Range("H110").Value = 7145.231
$ws.Range("I110").Value = 5888.9
$ws.Range("J110").Value = 11333
$ws.Range("K110").Value = 5888.9
$ws.Range("L110").Value = 11333
$ws.Range("M110").Value = -3843.9
$ws.Range("N110").Value = -15423
$ws.Range("H116").Value = 2850.077
$ws.Range("I116").Value = 2196.5789
$ws.Range("J116").Value = 4623.857
$ws.Range("K116").Value = 2196.5789
$ws.Range("L116").Value = 4623.857
$ws.Range("M116").Value = 97.42110000000002
$ws.Range("N116").Value = -9211.857
$ws.Range("H122").Value = 2736
$ws.Range("I122").Value = 2706.6667
$ws.Range("K122").Value = 8120.000100000001
$ws.Range("M122").Value = -5670.000100000001
$ws.Range("H132").Value = 4831.939
$ws.Range("I132").Value = 4585.75
$ws.Range("K132").Value = 13757.25
$ws.Range("M132").Value = -11227.25
$ws.Range("H136").Value = 19236778
$ws.Range("I136").Value = 22732012
$ws.Range("K136").Value = 68196036
$ws.Range("M136").Value = -68193486

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2850.077
$ws.Range("I3").Value = 2196.5789
$ws.Range("J3").Value = 4623.857
$ws.Range("K3").Value = 2196.5789
$ws.Range("L3").Value = 4623.857
$ws.Range("M3").Value = -2082.5789
$ws.Range("N3").Value = -4851.857
$ws.Range("H99").Value = 3285.4
$ws.Range("I99").Value = 2585.5557
$ws.Range("J99").Value = 4335.1665
$ws.Range("K99").Value = 2585.5557
$ws.Range("L99").Value = 4335.1665
$ws.Range("M99").Value = -1087.5557
$ws.Range("N99").Value = -7331.1665
$ws.Range("H132").Value = 98000
$ws.Range("J132").Value = 98000
$ws.Range("L132").Value = 98000
$ws.Range("N132").Value = -108120
$ws.Range("H134").Value = 2592
$ws.Range("I134").Value = 2575.68
$ws.Range("K134").Value = 7727.039999999999
$ws.Range("M134").Value = -5192.039999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3665.52
$ws.Range("I31").Value = 1697.3077
$ws.Range("J31").Value = 5797.75
$ws.Range("K31").Value = 1697.3077
$ws.Range("L31").Value = 5797.75
$ws.Range("M31").Value = -1402.3077
$ws.Range("N31").Value = -6387.75
$ws.Range("H34").Value = 3665.52
$ws.Range("I34").Value = 1697.3077
$ws.Range("J34").Value = 5797.75
$ws.Range("K34").Value = 1697.3077
$ws.Range("L34").Value = 5797.75
$ws.Range("M34").Value = -1495.3077
$ws.Range("N34").Value = -6201.75
$ws.Range("H97").Value = 33323.332
$ws.Range("I97").Value = 30000
$ws.Range("J97").Value = 34985
$ws.Range("K97").Value = 30000
$ws.Range("L97").Value = 34985
$ws.Range("M97").Value = -29009
$ws.Range("N97").Value = -36967
$ws.Range("H105").Value = 1910.75
$ws.Range("I105").Value = 2059.8
$ws.Range("K105").Value = 2059.8
$ws.Range("M105").Value = -312.8000000000002
$ws.Range("H107").Value = 423.125
$ws.Range("I107").Value = 302.13333
$ws.Range("J107").Value = 624.7778
$ws.Range("K107").Value = 302.13333
$ws.Range("L107").Value = 624.7778
$ws.Range("M107").Value = 1617.86667
$ws.Range("N107").Value = -4464.7778

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 21083774
$ws.Range("I4").Value = 21165262
$ws.Range("J4").Value = 20603000
$ws.Range("K4").Value = 63495786
$ws.Range("L4").Value = 61809000
$ws.Range("M4").Value = -63495674
$ws.Range("N4").Value = -61809224
$ws.Range("H39").Value = 555.5
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H107").Value = 1242.6364
$ws.Range("I107").Value = 1366.5
$ws.Range("K107").Value = 4099.5
$ws.Range("M107").Value = -2179.5
$ws.Range("H122").Value = 6219
$ws.Range("J122").Value = 1625.3334
$ws.Range("L122").Value = 14628.0006
$ws.Range("N122").Value = -19528.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 49997
$ws.Range("I53").Value = 49995
$ws.Range("K53").Value = 49995
$ws.Range("M53").Value = -49364
$ws.Range("H74").Value = 48000.6
$ws.Range("J74").Value = 48000.6
$ws.Range("L74").Value = 48000.6
$ws.Range("N74").Value = -49872.6
$ws.Range("H77").Value = 48000.6
$ws.Range("J77").Value = 48000.6
$ws.Range("L77").Value = 144001.8
$ws.Range("N77").Value = -153361.8
$ws.Range("H122").Value = 2849.762
$ws.Range("I122").Value = 2623.4211
$ws.Range("K122").Value = 7870.263300000001
$ws.Range("M122").Value = -5420.263300000001
$ws.Range("H126").Value = 4662.4165
$ws.Range("I126").Value = 3618.75
$ws.Range("J126").Value = 6749.75
$ws.Range("K126").Value = 10856.25
$ws.Range("L126").Value = 20249.25
$ws.Range("M126").Value = -8386.25
$ws.Range("N126").Value = -25189.25
$ws.Range("H132").Value = 5252.479
$ws.Range("I132").Value = 3725.72
$ws.Range("K132").Value = 11177.16
$ws.Range("M132").Value = -8647.16

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 21666.445
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312
$ws.Range("H62").Value = 209997.5
$ws.Range("J62").Value = 209997.5
$ws.Range("L62").Value = 209997.5
$ws.Range("N62").Value = -211245.5
$ws.Range("H65").Value = 209997.5
$ws.Range("J65").Value = 209997.5
$ws.Range("L65").Value = 629992.5
$ws.Range("N65").Value = -636232.5
$ws.Range("H82").Value = 8319.034
$ws.Range("I82").Value = 10148.875
$ws.Range("K82").Value = 10148.875
$ws.Range("M82").Value = -9787.875
$ws.Range("H85").Value = 8319.034
$ws.Range("I85").Value = 10148.875
$ws.Range("K85").Value = 10148.875
$ws.Range("M85").Value = -8900.875
$ws.Range("H122").Value = 2433.3333
$ws.Range("I122").Value = 2433.3333
$ws.Range("K122").Value = 7299.999899999999
$ws.Range("M122").Value = -4849.999899999999
$ws.Range("H132").Value = 18951
$ws.Range("I132").Value = 30719.8
$ws.Range("J132").Value = 7182.2
$ws.Range("K132").Value = 92159.39999999999
$ws.Range("L132").Value = 21546.6
$ws.Range("M132").Value = -89629.39999999999
$ws.Range("N132").Value = -26606.6
$ws.Range("H136").Value = 3450.7144
$ws.Range("I136").Value = 3217.5
$ws.Range("J136").Value = 4850
$ws.Range("K136").Value = 9652.5
$ws.Range("L136").Value = 14550
$ws.Range("M136").Value = -7102.5
$ws.Range("N136").Value = -19650

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 45000
$ws.Range("I45").Value = 45000
$ws.Range("K45").Value = 45000
$ws.Range("M45").Value = -44509
$ws.Range("H70").Value = 39111.25
$ws.Range("I70").Value = 34995
$ws.Range("J70").Value = 40483.332
$ws.Range("K70").Value = 34995
$ws.Range("L70").Value = 40483.332
$ws.Range("M70").Value = -34680
$ws.Range("N70").Value = -41113.332
$ws.Range("H73").Value = 39111.25
$ws.Range("I73").Value = 34995
$ws.Range("J73").Value = 40483.332
$ws.Range("K73").Value = 34995
$ws.Range("L73").Value = 40483.332
$ws.Range("M73").Value = -33903
$ws.Range("N73").Value = -42667.332
$ws.Range("H132").Value = 2595.5134
$ws.Range("I132").Value = 2436.6858
$ws.Range("J132").Value = 5375
$ws.Range("K132").Value = 7310.057400000001
$ws.Range("L132").Value = 16125
$ws.Range("M132").Value = -4780.057400000001
$ws.Range("N132").Value = -21185
